$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# ------------------------------------------------------------------
# 1. Insert a new blank row at row 21 (inside the leave table), which
#    shifts the existing rows 21-131 down to 22-132.
# ------------------------------------------------------------------
$ws.Rows(21).Insert() | Out-Null

# The inserted row does not inherit the surrounding table formatting,
# so copy the number formats / borders from the row below (row 22,
# which is the original row 21 that just got pushed down).
$ws.Range("A22:K22").Copy() | Out-Null
$ws.Range("A21:K21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Restore the calculated-column formula in G21 (lost by the formats-only paste)
$ws.Range("G21").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# ------------------------------------------------------------------
# 2. Grow the Table1 list object so it covers the new row at the end
#    (A8:K131 -> A8:K132), mirroring what Excel does automatically
#    when a row is inserted in the middle of a table.
# ------------------------------------------------------------------
$lo.Resize($ws.Range("A8:K132")) | Out-Null

# Restore the calculated-column formula text in the brand new last row
# (the Resize operation fills it in with a different, non-evaluating
# syntax) so it matches the rest of the column.
$ws.Range("G132").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# Row 130 keeps the regular (non-last-row) styling, row 131 only
# changes its E/I styling to the regular style, and the brand new
# row 132 becomes the new "last row" of the table, carrying the
# special border styles that used to belong to row 131.
$ws.Range("E130").Copy() | Out-Null
$ws.Range("E131").PasteSpecial(-4122) | Out-Null
$ws.Range("I130").Copy() | Out-Null
$ws.Range("I131").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Fill in the new leave entry itself.
# ------------------------------------------------------------------
$ws.Range("B21").Value = "FL(2-0-0)"
$ws.Range("D21").Value = 2
$ws.Range("K21").Value = "10/13,16/2023"

# Row 19 (9/1/2023 period) now has an EARNED value of 1.25 as well.
$ws.Range("C19").Value = 1.25

# ------------------------------------------------------------------
# 4. Recalculate so all the dependent formulas (BALANCE columns in
#    Sheet1 row 9 and CONVERTION!A7) pick up the new totals.
# ------------------------------------------------------------------
$excel.CalculateFullRebuild() | Out-Null

# Leave the selection on the newly added remark cell, as the author did.
$ws.Activate() | Out-Null
$ws.Range("K21").Select() | Out-Null
